$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Clear the "J" status column (Pass/Fail values) for data rows 2-17,
# leaving the "Results" header in J1 intact.
$ws.Range("J2:J17").ClearContents()

# Reset the view/selection to the top of the sheet.
$ws.Activate()
$ws.Range("A2").Select()
